$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 20
$ws.Range("A20").Value = "lav UI til UC01"
$ws.Range("C20").Value = 43893
$ws.Range("D20").Value = 0.36805555555555558
$ws.Range("E20").Value = 0.41666666666666669

# Row 21
$ws.Range("A21").Value = "lav UI til UC08"
$ws.Range("C21").Value = 43893
$ws.Range("D21").Value = 0.57638888888888895
$ws.Range("E21").Value = 0.61111111111111105

# Row 22
$ws.Range("A22").Value = "undersøgt maven mulighed"
$ws.Range("B22").Value = "Software Architect"
$ws.Range("C22").Value = 43893
$ws.Range("D22").Value = 0.62847222222222221
$ws.Range("E22").Value = 0.66319444444444442

# Update sheet view: scroll position and selection
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("E23").Select()
